# In vitro validation 1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-13 (columns A-K)
$data = @(
    @(12292, "18PG",  0, 9.3028,             10,  11,  20,  55,  4,   -0.03219976809333189, 1.803581414630675),
    @(13512, "18PG",  0, 9.3028,             10,  11,  40,  35,  70,  -0.02076262833375528, 1.873347967164093),
    @(2711,  "18PG",  0, 9.3028,             10,  6,   200, 50,  4,   -0.01909817385362694, 1.883501139492876),
    @(7506,  "18PG",  0, 9.3028,             10,  9,   40,  55,  4,   -0.005086681177133375, 1.968971244819486),
    @(28668, "DOPE",  1, 9.784000000000001,  9,   11,  250, 50,  180, 0.9066252696552797,  7.530414144897206),
    @(23740, "DOPE",  1, 9.784000000000001,  9,   9,   250, 45,  160, 0.9833729197174315,  7.998574810276332),
    @(14765, "DOPE",  1, 9.784000000000001,  9,   5,   200, 30,  5,   0.9135837784148289,  7.572861048330456),
    @(26131, "DOPE",  1, 9.784000000000001,  9,   10,  90,  45,  160, 0.9833729197174315,  7.998574810276332),
    @(26186, "DOPE",  1, 9.784000000000001,  9,   10,  80,  70,  160, 0.8790140722709167,  7.361985840852591),
    @(14608, "DOPE",  1, 9.784000000000001,  9,   5,   80,  75,  3,   0.8799861653012702,  7.367915608337746),
    @(49643, "DSPC",  1, 6.12,               8,   7,   450, 75,  80,  -0.02254784829520736, 1.862458125399235),
    @(47698, "DSPC",  1, 6.12,               8,   6,   20,  55,  160, -0.02821529602436637, 1.827886694251365)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le 11; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
